# Apply the "add MoransI and fix codes" edits to methods.xlsx
$wb = $excel.ActiveWorkbook
$wsM = $wb.Worksheets.Item("Methods")
$wsP = $wb.Worksheets.Item("Parameters")

# ---------------------------------------------------------------
# 1. Fix mis-copied "Model" / "Model Category" values (NB/Poisson -> GP,
#    probabilistic model -> GAMLSS, ZINB/Optimisation framework -> Normal
#    distribution, and dyntoy's NB -> ZINB / Zero-inflated model).
# ---------------------------------------------------------------

# powsimR (row 9): NB, NB -> GP, GP
$wsM.Range("E9").Value2 = "GP"
$wsM.Range("F9").Value2 = "GP"

# SCRIP-paths (row 14): Poisson, Poisson -> GP, GP
$wsM.Range("E14").Value2 = "GP"
$wsM.Range("F14").Value2 = "GP"

# scDesign3-tree (row 18): probabilistic model / Probabilistic model -> GAMLSS
$wsM.Range("E18").Value2 = "GAMLSS"
$wsM.Range("F18").Value2 = "GAMLSS"

# PROSSTT (row 22): NB, NB -> GP, GP
$wsM.Range("E22").Value2 = "GP"
$wsM.Range("F22").Value2 = "GP"

# dyntoy (row 23): NB, NB -> ZINB, Zero-inflated model
$wsM.Range("E23").Value2 = "ZINB"
$wsM.Range("F23").Value2 = "Zero-inflated model"

# MFA (row 25): ZINB, Zero-inflated model -> Normal distribution, Other models
$wsM.Range("E25").Value2 = "Normal distribution"
$wsM.Range("F25").Value2 = "Other models"

# phenopath (row 26): NB, NB -> GP, GP
$wsM.Range("E26").Value2 = "GP"
$wsM.Range("F26").Value2 = "GP"

# Lun (row 27): NB, NB -> GP, GP
$wsM.Range("E27").Value2 = "GP"
$wsM.Range("F27").Value2 = "GP"

# scDD (row 28): NB, NB -> GP, GP
$wsM.Range("E28").Value2 = "GP"
$wsM.Range("F28").Value2 = "GP"

# muscat (row 31): NB, NB -> GP, GP
$wsM.Range("E31").Value2 = "GP"
$wsM.Range("F31").Value2 = "GP"

# scDesign3 (row 34): probabilistic model / Probabilistic model -> GAMLSS
$wsM.Range("E34").Value2 = "GAMLSS"
$wsM.Range("F34").Value2 = "GAMLSS"

# SparseDC (row 38): Optimisation framework -> Normal distribution (F38 unchanged: Other models)
$wsM.Range("E38").Value2 = "Normal distribution"

# hierarchicell (row 39): NB, NB -> GP, GP
$wsM.Range("E39").Value2 = "GP"
$wsM.Range("F39").Value2 = "GP"

# Lun2 (row 43): NB, NB -> GP, GP
$wsM.Range("E43").Value2 = "GP"
$wsM.Range("F43").Value2 = "GP"

# BASiCS (row 44): NB, NB -> GP, GP
$wsM.Range("E44").Value2 = "GP"
$wsM.Range("F44").Value2 = "GP"

# Simple (row 46): NB, NB -> GP, GP
$wsM.Range("E46").Value2 = "GP"
$wsM.Range("F46").Value2 = "GP"

# Kersplat (row 47): Poisson, Poisson -> GP, GP
$wsM.Range("E47").Value2 = "GP"
$wsM.Range("F47").Value2 = "GP"

# ---------------------------------------------------------------
# 2. MFA row (25) gets a checkmark for "Simulate Trajectory" (K25) using a
#    distinct symbol font (Segoe UI Symbol, 14pt) and without wrap text,
#    and the row grows a bit taller to fit it.
# ---------------------------------------------------------------
$wsM.Range("K25").Value2 = "$([char]0x2713)"
$wsM.Range("K25").Font.Name = "Segoe UI Symbol"
$wsM.Range("K25").Font.Size = 14
$wsM.Range("K25").WrapText = $false
$wsM.Range("K25").HorizontalAlignment = -4108
$wsM.Range("K25").VerticalAlignment = -4108
$wsM.Rows.Item(25).RowHeight = 21

# scDD row (28) loses its special wrap/symbol font on E28:F28 (now regular body style)
$wsM.Range("E28:F28").WrapText = $false
$wsM.Range("E28:F28").Font.Name = "Times New Roman"
$wsM.Range("E28:F28").Font.Size = 12
$wsM.Range("E28:F28").HorizontalAlignment = -4108
$wsM.Range("E28:F28").VerticalAlignment = -4108
$wsM.Rows.Item(28).RowHeight = 18

# ---------------------------------------------------------------
# 3. Notes section: remove the "NB: Negative Binomial" row and fix the
#    ZINB note's wording (zero-inflated normal -> zero-inflated negative
#    binomial). This shifts rows 58-61 up to 57-60.
# ---------------------------------------------------------------
$wsM.Rows.Item(57).Delete()
$wsM.Range("A57").Value2 = "ZINB: zero-inflated negative binomial distribution"

# ---------------------------------------------------------------
# 4. View / selection state: Methods becomes the active sheet/tab, with a
#    new selection; Parameters keeps a selection but is no longer the
#    active tab.
# ---------------------------------------------------------------
$wsP.Activate()
$wsP.Range("D34").Select()

$wsM.Activate()
$wsM.Range("E15").Select()
